$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "sides_arrowsmith"

# --- Write the full data grid ---
$ws.Cells.Item(1, 1).Value = "ItemName"
$ws.Cells.Item(1, 2).Value = "Ingredients"
$ws.Cells.Item(1, 3).Value = "Allergens"
$ws.Cells.Item(1, 4).Value = "LocalIngredients"
$ws.Cells.Item(1, 5).Value = "Diet"
$ws.Cells.Item(1, 6).Value = "Nutrition Label"
$ws.Cells.Item(1, 7).Value = "LeaveEmpty"

$ws.Cells.Item(2, 1).Value = "Fries"
$ws.Cells.Item(2, 2).Value = "McCain Gold Crisp Fries[Potatoes, canola oil, wheat flour, modified corn starch, corn flour, salt, autolyzed yeast, baking powder, sodium phosphate, dextrose, modified cellulose, colour]"
$ws.Cells.Item(2, 3).Value = "Wheat."
$ws.Cells.Item(2, 4).Value = "NA"
$ws.Cells.Item(2, 5).Value = "VEG,VGN,DF"
$ws.Cells.Item(2, 6).Value = "Fries"

$ws.Cells.Item(3, 1).Value = "Poutine"
$ws.Cells.Item(3, 2).Value = " Fries / Gravy / Cheese Curds"
$ws.Cells.Item(3, 3).Value = "Wheat, milk, soy."
$ws.Cells.Item(3, 4).Value = "NA"
$ws.Cells.Item(3, 5).Value = "VEG,DF"
$ws.Cells.Item(3, 6).Value = "placeholder"

$ws.Cells.Item(4, 1).Value = "Extra Cheese Curds"
$ws.Cells.Item(4, 2).Value = "29% MF Cheese Curds"
$ws.Cells.Item(4, 3).Value = "Milk."
$ws.Cells.Item(4, 4).Value = "NA"
$ws.Cells.Item(4, 5).Value = "VEG"
$ws.Cells.Item(4, 6).Value = "placeholder"

$ws.Cells.Item(5, 1).Value = "Side Gravy"
$ws.Cells.Item(5, 2).Value = "Vegetarian Brown Gravy"
$ws.Cells.Item(5, 3).Value = "Wheat, milk, soy."
$ws.Cells.Item(5, 4).Value = "NA"
$ws.Cells.Item(5, 5).Value = "VEG"
$ws.Cells.Item(5, 6).Value = "placeholder"

$ws.Cells.Item(6, 1).Value = "Tater Tot"
$ws.Cells.Item(6, 2).Value = "Potatoes / Canola Oil / Salt / Seasonings"
$ws.Cells.Item(6, 3).Value = "No known priotity allergens."
$ws.Cells.Item(6, 4).Value = "NA"
$ws.Cells.Item(6, 5).Value = "VEG,VGN,GF,DF"
$ws.Cells.Item(6, 6).Value = "Tater_Tots"

$ws.Cells.Item(7, 1).Value = "12 oz Soup"
$ws.Cells.Item(7, 2).Value = "Item will vary daily - find information in the cafeteria."
$ws.Cells.Item(7, 3).Value = "Item will vary daily - find information in the cafeteria."
$ws.Cells.Item(7, 4).Value = "NA"
$ws.Cells.Item(7, 5).Value = "NA"
$ws.Cells.Item(7, 6).Value = "placeholder"

# --- Apply the "wrap text" style (existing style, reused) to the Ingredients column ---
$ws.Range("B3,B4,B5").WrapText = $true

# --- Apply the new centered/bold-ish style to the ItemName column on certain rows ---
$ws.Range("A2,A6,A7").HorizontalAlignment = -4108
$ws.Range("A2,A6,A7").VerticalAlignment = -4108
$ws.Range("A2,A6,A7").Font.Size = 11
$ws.Range("A2,A6,A7").Font.Color = 0

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()

# --- Resize the table to cover the new data range ---
$t = $ws.ListObjects.Item(1)
$t.Resize($ws.Range("A1:G7"))

# --- Update selection to match the final saved state ---
$ws.Range("B18").Select()

